$d = $word.ActiveDocument
$sec = $d.Sections(1)

# Header "first page" (header1.xml) holds the BTec_Logo-Orange picture,
# currently named "image2.jpg" -> rename to "image1.jpg".
$hdrFirst = $sec.Headers(2)
if ($hdrFirst.Exists -and $hdrFirst.Range.InlineShapes.Count -ge 1) {
    $btecLogo = $hdrFirst.Range.InlineShapes.Item(1)
    $btecLogo.Name = "image1.jpg"
}

# Footer "default" (footer2.xml) holds a PearsonLogo.png picture (docPr id="2"),
# currently named "image1.png" -> rename to "image2.png".
$ftrDefault = $sec.Footers(1)
if ($ftrDefault.Exists -and $ftrDefault.Range.InlineShapes.Count -ge 1) {
    $pearsonLogoDefault = $ftrDefault.Range.InlineShapes.Item(1)
    $pearsonLogoDefault.Name = "image2.png"
}

# Footer "first page" (footer1.xml) holds a PearsonLogo.png picture (docPr id="3"),
# currently named "image1.png" -> rename to "image2.png".
$ftrFirst = $sec.Footers(2)
if ($ftrFirst.Exists -and $ftrFirst.Range.InlineShapes.Count -ge 1) {
    $pearsonLogoFirst = $ftrFirst.Range.InlineShapes.Item(1)
    $pearsonLogoFirst.Name = "image2.png"
}
